$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -1226

# Row 38
$ws.Range("H38").Value = 630.1053000000001
$ws.Range("I38").Value = 94.07692
$ws.Range("J38").Value = 1791.5
$ws.Range("K38").Value = 282.23076
$ws.Range("L38").Value = 5374.5
$ws.Range("M38").Value = 89.76923999999997
$ws.Range("N38").Value = -6118.5

# Row 64
$ws.Range("H64").Value = 3374.8
$ws.Range("I64").Value = 3181.4546
$ws.Range("J64").Value = 3611.111
$ws.Range("K64").Value = 3181.4546
$ws.Range("L64").Value = 3611.111
$ws.Range("M64").Value = -2933.4546
$ws.Range("N64").Value = -4107.111

# Row 67
$ws.Range("H67").Value = 3374.8
$ws.Range("I67").Value = 3181.4546
$ws.Range("J67").Value = 3611.111
$ws.Range("K67").Value = 3181.4546
$ws.Range("L67").Value = 3611.111
$ws.Range("M67").Value = -2323.4546
$ws.Range("N67").Value = -5327.111

# Row 76
$ws.Range("H76").Value = 2966.0232
$ws.Range("I76").Value = 2616.5
$ws.Range("J76").Value = 3022.7026
$ws.Range("K76").Value = 2616.5
$ws.Range("L76").Value = 3022.7026
$ws.Range("M76").Value = -2301.5
$ws.Range("N76").Value = -3652.7026

# Row 79
$ws.Range("H79").Value = 2966.0232
$ws.Range("I79").Value = 2616.5
$ws.Range("J79").Value = 3022.7026
$ws.Range("K79").Value = 2616.5
$ws.Range("L79").Value = 3022.7026
$ws.Range("M79").Value = -1524.5
$ws.Range("N79").Value = -5206.702600000001

# Row 124
$ws.Range("H124").Value = 33903.8
$ws.Range("J124").Value = 33903.8
$ws.Range("L124").Value = 33903.8
$ws.Range("N124").Value = -43723.8

# Row 125
$ws.Range("H125").Value = 695.625
$ws.Range("I125").Value = 738.5
$ws.Range("J125").Value = 652.75
$ws.Range("K125").Value = 6646.5
$ws.Range("L125").Value = 5874.75
$ws.Range("M125").Value = -4186.5
$ws.Range("N125").Value = -10794.75

# Row 126
$ws.Range("H126").Value = 79675
$ws.Range("J126").Value = 79675
$ws.Range("L126").Value = 79675
$ws.Range("N126").Value = -89555

# Row 127
$ws.Range("H127").Value = 889.4516
$ws.Range("I127").Value = 497.8
$ws.Range("J127").Value = 1075.9524
$ws.Range("K127").Value = 1493.4
$ws.Range("L127").Value = 3227.857199999999
$ws.Range("M127").Value = 3466.6
$ws.Range("N127").Value = -13147.8572

# Row 130
$ws.Range("H130").Value = 51780
$ws.Range("J130").Value = 51780
$ws.Range("L130").Value = 51780
$ws.Range("N130").Value = -61820

# Row 133
$ws.Range("H133").Value = 38312.5
$ws.Range("J133").Value = 38312.5
$ws.Range("L133").Value = 38312.5
$ws.Range("N133").Value = -48432.5

# Row 139
$ws.Range("H139").Value = 11111
$ws.Range("J139").Value = 11111
$ws.Range("L139").Value = 11111
$ws.Range("N139").Value = -21391

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 11427.571
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 11427.571
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 11427.571
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -11773.571

# Row 63
$ws.Range("H63").Value = 7000
$ws.Range("I63").Value = 5000
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 5000
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -4314
$ws.Range("N63").Value = -9372

# Row 66
$ws.Range("H66").Value = 7000
$ws.Range("I66").Value = 5000
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 25000
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -21568
$ws.Range("N66").Value = -46864

# Row 69
$ws.Range("H69").Value = 43000
$ws.Range("J69").Value = 43000
$ws.Range("L69").Value = 43000
$ws.Range("N69").Value = -44498

# Row 72
$ws.Range("H72").Value = 43000
$ws.Range("J72").Value = 43000
$ws.Range("L72").Value = 129000
$ws.Range("N72").Value = -136488

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 10676.454
$ws.Range("I26").Value = 7493.4443
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 7493.4443
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = -7201.4443
$ws.Range("N26").Value = -25584

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 5683.143
$ws.Range("I3").Value = 945.5
$ws.Range("K3").Value = 945.5
$ws.Range("M3").Value = -832.5

$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 3722.0715
$ws.Range("I63").Value = 3447.5
$ws.Range("J63").Value = 3767.8333
$ws.Range("K63").Value = 10342.5
$ws.Range("L63").Value = 11303.4999
$ws.Range("M63").Value = -9593.5
$ws.Range("N63").Value = -12801.4999

# Row 66
$ws.Range("H66").Value = 3722.0715
$ws.Range("I66").Value = 3447.5
$ws.Range("J66").Value = 3767.8333
$ws.Range("K66").Value = 31027.5
$ws.Range("L66").Value = 33910.4997
$ws.Range("M66").Value = -27283.5
$ws.Range("N66").Value = -41398.4997

# Row 68
$ws.Range("H68").Value = 1158.6586
$ws.Range("I68").Value = 1106
$ws.Range("K68").Value = 3318
$ws.Range("M68").Value = -2507

# Row 71
$ws.Range("H71").Value = 1158.6586
$ws.Range("I71").Value = 1106
$ws.Range("K71").Value = 9954
$ws.Range("M71").Value = -5898

# Row 87
$ws.Range("H87").Value = 28368.809
$ws.Range("I87").Value = 1902.4
$ws.Range("J87").Value = 29776.596
$ws.Range("K87").Value = 5707.200000000001
$ws.Range("L87").Value = 89329.788
$ws.Range("M87").Value = -4459.200000000001
$ws.Range("N87").Value = -91825.788

# Row 90
$ws.Range("H90").Value = 28368.809
$ws.Range("I90").Value = 1902.4
$ws.Range("J90").Value = 29776.596
$ws.Range("K90").Value = 17121.6
$ws.Range("L90").Value = 267989.364
$ws.Range("M90").Value = -10881.6
$ws.Range("N90").Value = -280469.364

# Row 93
$ws.Range("H93").Value = 4804.952
$ws.Range("I93").Value = 3612
$ws.Range("J93").Value = 4930.5264
$ws.Range("K93").Value = 10836
$ws.Range("L93").Value = 14791.5792
$ws.Range("N93").Value = -18535.5792
$ws.Range("M93").Value = -8964

# Row 131
$ws.Range("H131").Value = 2532.4033
$ws.Range("I131").Value = 567.5
$ws.Range("J131").Value = 2667.9138
$ws.Range("K131").Value = 1702.5
$ws.Range("L131").Value = 8003.741399999999
$ws.Range("M131").Value = 3337.5
$ws.Range("N131").Value = -18083.7414

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 166668180
$ws.Range("I46").Value = 250000770
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 250000770
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -250000582
$ws.Range("N46").Value = -3376

# Row 74
$ws.Range("H74").Value = 18000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 18000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 18000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -19996

# Row 77
$ws.Range("H77").Value = 18000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 18000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 54000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -63984

# Row 132
$ws.Range("H132").Value = 5111.425
$ws.Range("I132").Value = 4980.8887
$ws.Range("J132").Value = 5382.5386
$ws.Range("K132").Value = 14942.6661
$ws.Range("L132").Value = 16147.6158
$ws.Range("M132").Value = -12412.6661
$ws.Range("N132").Value = -21207.6158
